# Remove the last slide ("글상자") from the presentation.
$p = $ppt.ActivePresentation
$count = $p.Slides.Count
$p.Slides.Item($count).Delete()
